$wb = $excel.ActiveWorkbook

# --- Sheet "2016" (sheet4.xml): new data row for Amazonas (AM), existing rows shift down ---
$ws2016 = $wb.Worksheets.Item("2016")

# Row 5 (AM) becomes brand-new data
$ws2016.Range("C5").Value = 115455
$ws2016.Range("E5").Value = 53710
$ws2016.Range("G5").Value = 319067
$ws2016.Range("I5").Value = 90039.7

# Row 6 (BA) takes what used to be row 5's data (rounded)
$ws2016.Range("C6").Value = 391881
$ws2016.Range("E6").Value = 225991
$ws2016.Range("G6").Value = 1196220
$ws2016.Range("I6").Value = 351638

# Row 7 (CE) takes what used to be row 6's data (rounded)
$ws2016.Range("C7").Value = 160089
$ws2016.Range("E7").Value = 67007
$ws2016.Range("G7").Value = 749490
$ws2016.Range("I7").Value = 201037

# Row 8 (DF) takes what used to be row 7's data (rounded)
$ws2016.Range("C8").Value = 78188
$ws2016.Range("E8").Value = 35876.300000000003
$ws2016.Range("G8").Value = 225223
$ws2016.Range("I8").Value = 61685.7

# Row 9 (ES) takes what used to be row 8's data (unchanged values)
$ws2016.Range("C9").Value = 89265.8
$ws2016.Range("E9").Value = 34710.400000000001
$ws2016.Range("G9").Value = 332085.40000000002
$ws2016.Range("I9").Value = 89914.6

# --- Active sheet / selection change: "2019" was active, now "2016" is active ---
$ws2016.Activate()
$ws2016.Range("F9").Select()
